$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.590.10"
$ws.Range("D3").Value = "1.680.32"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.525"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.06"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0627"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "1.921.41"
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +12.67%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.621"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +9.44%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.690.31"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.98"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "30.577.49"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "0.0₃0716"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.28"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.79"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.112"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.48"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("D33").Value = "1.508.79"
$ws.Range("E33").Value = "  +5.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.29"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.75"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "84.08"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.43%  "
$ws.Range("E38").Value = "  +4.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.588"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.838"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.99"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0500"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.02"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "51.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.23%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.51"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.814.26"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("E50").Value = "  +5.61%  "
$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +1.94%  "
